$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-09 03:48:23"
$ws.Range("N2").Value = "-5.3 °C 3:05 TU"
$ws.Range("O2").Value = "-3.8 °C"
$ws.Range("E3").Value = "2026-02-09 03:48:26"
$ws.Range("O3").Value = "-6.5 °C"
$ws.Range("E4").Value = "2026-02-09 03:48:28"
$ws.Range("N4").Value = "3.6 °C 3:29 TU"
$ws.Range("O4").Value = "4.8 °C"
$ws.Range("E5").Value = "2026-02-09 03:48:30"
$ws.Range("N5").Value = "-6.3 °C 3:10 TU"
$ws.Range("O5").Value = "-5.5 °C"
$ws.Range("E6").Value = "2026-02-09 03:48:33"
$ws.Range("N6").Value = "5.8 °C 3:29 TU"
$ws.Range("O6").Value = "7.0 °C"
$ws.Range("E7").Value = "2026-02-09 03:48:35"
$ws.Range("E8").Value = "2026-02-09 03:48:38"
$ws.Range("H8").Value = "'82%"
$ws.Range("E9").Value = "2026-02-09 03:48:40"
$ws.Range("H9").Value = "'81%"
$ws.Range("N9").Value = "5.0 °C 3:29 TU"
$ws.Range("O9").Value = "8.1 °C"
$ws.Range("E10").Value = "2026-02-09 03:48:43"
$ws.Range("H10").Value = "'93%"
$ws.Range("N10").Value = "3.8 °C 3:13 TU"
$ws.Range("O10").Value = "6.5 °C"
$ws.Range("E11").Value = "2026-02-09 03:48:45"
$ws.Range("M11").Value = "2.7 °C 3:29 TU"
$ws.Range("O11").Value = "1.8 °C"
$ws.Range("E12").Value = "2026-02-09 03:48:48"
$ws.Range("H12").Value = "'86%"
$ws.Range("N12").Value = "6.1 °C 3:28 TU"
$ws.Range("O12").Value = "8.2 °C"
$ws.Range("E13").Value = "2026-02-09 03:48:50"
$ws.Range("N13").Value = "-2.7 °C 3:09 TU"
$ws.Range("O13").Value = "-1.2 °C"
$ws.Range("E14").Value = "2026-02-09 03:48:52"
$ws.Range("H14").Value = "'97%"
$ws.Range("E15").Value = "2026-02-09 03:48:55"
$ws.Range("H15").Value = "'84%"
$ws.Range("N15").Value = "4.1 °C 3:29 TU"
$ws.Range("O15").Value = "6.5 °C"
$ws.Range("E16").Value = "2026-02-09 03:48:57"
$ws.Range("H16").Value = "'63%"
$ws.Range("L16").Value = "20.2 km/h - 229º 3:22 TU"
$ws.Range("M16").Value = "-4.5 °C 3:28 TU"
$ws.Range("N16").Value = "-6.1 °C 3:00 TU"
$ws.Range("O16").Value = "-5.2 °C"
$ws.Range("E17").Value = "2026-02-09 03:48:59"
$ws.Range("H17").Value = "'93%"
$ws.Range("N17").Value = "-0.7 °C 3:22 TU"
$ws.Range("O17").Value = "0.1 °C"
$ws.Range("E18").Value = "2026-02-09 03:49:02"
$ws.Range("N18").Value = "5.8 °C 3:02 TU"
$ws.Range("O18").Value = "7.0 °C"
$ws.Range("E19").Value = "2026-02-09 03:49:04"
$ws.Range("E20").Value = "2026-02-09 03:49:07"
$ws.Range("N20").Value = "-7.7 °C 3:29 TU"
$ws.Range("O20").Value = "-6.2 °C"
$ws.Range("E21").Value = "2026-02-09 03:49:10"
$ws.Range("H21").Value = "'92%"
$ws.Range("N21").Value = "-0.4 °C 3:29 TU"
$ws.Range("O21").Value = "1.0 °C"
$ws.Range("E22").Value = "2026-02-09 03:49:12"
$ws.Range("N22").Value = "-8.0 °C 3:12 TU"
$ws.Range("O22").Value = "-7.7 °C"
$ws.Range("E23").Value = "2026-02-09 03:49:14"
$ws.Range("N23").Value = "-6.9 °C 3:26 TU"
$ws.Range("O23").Value = "-5.8 °C"
$ws.Range("E24").Value = "2026-02-09 03:49:17"
$ws.Range("H24").Value = "'83%"
$ws.Range("N24").Value = "3.1 °C 3:10 TU"
$ws.Range("O24").Value = "4.8 °C"
$ws.Range("E25").Value = "2026-02-09 03:49:19"
$ws.Range("H25").Value = "'79%"
$ws.Range("M25").Value = "-3.5 °C 3:25 TU"
$ws.Range("O25").Value = "-4.6 °C"
$ws.Range("E26").Value = "2026-02-09 03:49:22"
$ws.Range("H26").Value = "'93%"
$ws.Range("E27").Value = "2026-02-09 03:49:24"
$ws.Range("H27").Value = "'83%"
$ws.Range("N27").Value = "-4.4 °C 3:16 TU"
$ws.Range("E28").Value = "2026-02-09 03:49:26"
$ws.Range("N28").Value = "2.5 °C 3:29 TU"
$ws.Range("O28").Value = "4.1 °C"
$ws.Range("E29").Value = "2026-02-09 03:49:29"
$ws.Range("O29").Value = "6.3 °C"
$ws.Range("E30").Value = "2026-02-09 03:49:32"
$ws.Range("H30").Value = "'93%"
$ws.Range("N30").Value = "5.8 °C 3:29 TU"
$ws.Range("O30").Value = "7.2 °C"
$ws.Range("E31").Value = "2026-02-09 03:49:34"
$ws.Range("H31").Value = "'74%"
$ws.Range("M31").Value = "9.7 °C 3:26 TU"
$ws.Range("O31").Value = "8.7 °C"
$ws.Range("E32").Value = "2026-02-09 03:49:37"
$ws.Range("H32").Value = "'78%"
$ws.Range("N32").Value = "2.7 °C 3:27 TU"
$ws.Range("E33").Value = "2026-02-09 03:49:39"
$ws.Range("H33").Value = "'95%"
$ws.Range("N33").Value = "-1.4 °C 3:14 TU"
$ws.Range("O33").Value = "-0.2 °C"
$ws.Range("E34").Value = "2026-02-09 03:49:41"
$ws.Range("H34").Value = "'74%"
$ws.Range("O34").Value = "-2.6 °C"
$ws.Range("E35").Value = "2026-02-09 03:49:44"
$ws.Range("H35").Value = "'68%"
$ws.Range("M35").Value = "4.1 °C 3:29 TU"
$ws.Range("E36").Value = "2026-02-09 03:49:47"
$ws.Range("H36").Value = "'77%"
$ws.Range("N36").Value = "6.1 °C 3:21 TU"
$ws.Range("O36").Value = "9.2 °C"
$ws.Range("E37").Value = "2026-02-09 03:49:49"
$ws.Range("H37").Value = "'92%"
$ws.Range("N37").Value = "3.1 °C 3:27 TU"
$ws.Range("O37").Value = "3.9 °C"
$ws.Range("E38").Value = "2026-02-09 03:49:51"
$ws.Range("H38").Value = "'98%"
$ws.Range("M38").Value = "7.3 °C 3:03 TU"
$ws.Range("O38").Value = "6.4 °C"
$ws.Range("E39").Value = "2026-02-09 03:49:54"
$ws.Range("H39").Value = "'86%"
$ws.Range("I39").Value = "0.1 mm"
$ws.Range("N39").Value = "-6.3 °C 3:09 TU"
$ws.Range("O39").Value = "-5.3 °C"
$ws.Range("E40").Value = "2026-02-09 03:49:56"
$ws.Range("N40").Value = "-1.1 °C 3:12 TU"
$ws.Range("O40").Value = "-0.1 °C"
$ws.Range("E41").Value = "2026-02-09 03:49:59"
$ws.Range("L41").Value = "22.0 km/h - 256º 3:12 TU"
$ws.Range("E42").Value = "2026-02-09 03:50:01"
$ws.Range("N42").Value = "5.5 °C 3:29 TU"
$ws.Range("O42").Value = "7.0 °C"
$ws.Range("E43").Value = "2026-02-09 03:50:04"
$ws.Range("E44").Value = "2026-02-09 03:50:06"
$ws.Range("N44").Value = "-9.0 °C 3:12 TU"
$ws.Range("O44").Value = "-7.7 °C"
$ws.Range("E45").Value = "2026-02-09 03:50:09"
$ws.Range("H45").Value = "'95%"
$ws.Range("J45").Value = "1010.0 hPa"
$ws.Range("L45").Value = "9.4 km/h - 140º 3:23 TU"
$ws.Range("E46").Value = "2026-02-09 03:50:11"
$ws.Range("N46").Value = "4.9 °C 3:01 TU"
$ws.Range("O46").Value = "6.1 °C"
